$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.133.93"

$ws.Range("D3").Value = "'1.848.53"
$ws.Range("E3").Value = "  -2.34%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'0.7024"
$ws.Range("E5").Value = "  -4.92%  "

$ws.Range("D6").Value = "'238.00"
$ws.Range("E6").Value = "  -1.96%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'0.3034"
$ws.Range("E8").Value = "  -4.38%  "

$ws.Range("D9").Value = "'0.07467"
$ws.Range("E9").Value = "  +3.52%  "

$ws.Range("D10").Value = "'23.34"
$ws.Range("E10").Value = "  -6.43%  "

$ws.Range("D11").Value = "'0.08123"
$ws.Range("E11").Value = "  -2.66%  "

$ws.Range("D12").Value = "'1.863.47"
$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").Value = "'0.7241"
$ws.Range("E13").Value = "  -4.80%  "

$ws.Range("D14").Value = "'5.203"
$ws.Range("E14").Value = "  -4.66%  "

$ws.Range("D15").Value = "'88.71"
$ws.Range("E15").Value = "  -4.74%  "

$ws.Range("D16").Value = "'29.114.81"
$ws.Range("E16").Value = "  -3.70%  "

$ws.Range("D17").Value = "'5.751"
$ws.Range("E17").Value = "  -6.87%  "

$ws.Range("D18").Value = "'236.93"
$ws.Range("E18").Value = "  -5.42%  "

$ws.Range("D19").Value = "'13.06"
$ws.Range("E19").Value = "  -4.29%  "

$ws.Range("D20").Value = "'0.000007638"
$ws.Range("E20").Value = "  -3.31%  "

$ws.Range("D21").Value = "'1.0000"
$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").Value = "'2.089.32"
$ws.Range("E22").Value = "  -4.52%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "'7.547"
$ws.Range("E24").Value = "  -5.12%  "

$ws.Range("D25").Value = "'8.973"
$ws.Range("E25").Value = "  -3.66%  "

$ws.Range("D26").Value = "'161.07"
$ws.Range("E26").Value = "  -2.20%  "

$ws.Range("D27").Value = "'0.1449"
$ws.Range("E27").Value = "  -8.57%  "

$ws.Range("D28").Value = "'18.01"
$ws.Range("E28").Value = "  -4.07%  "

$ws.Range("D29").Value = "'1.956"
$ws.Range("E29").Value = "  -5.22%  "

$ws.Range("D30").Value = "'1.396"
$ws.Range("E30").Value = "  -6.05%  "

$ws.Range("D31").Value = "'4.505"
$ws.Range("E31").Value = "  -1.87%  "

$ws.Range("D32").Value = "'1.491"
$ws.Range("E32").Value = "  -3.05%  "

$ws.Range("D33").Value = "'3.960"
$ws.Range("E33").Value = "  -5.89%  "

$ws.Range("D34").Value = "'0.05140"
$ws.Range("E34").Value = "  -4.32%  "

$ws.Range("D35").Value = "'1.183"
$ws.Range("E35").Value = "  -5.96%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6976"
$ws.Range("E36").Value = "  -10.52%  "

$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.013"
$ws.Range("E37").Value = "  +0.70%  "

$ws.Range("D38").Value = "'2.656"
$ws.Range("E38").Value = "  -2.82%  "

$ws.Range("D39").Value = "'0.01867"
$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("D40").Value = "'2.681"
$ws.Range("E40").Value = "  -3.04%  "

$ws.Range("D41").Value = "'0.9415"
$ws.Range("E41").Value = "  +7.87%  "

$ws.Range("D42").Value = "'5.995"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("D43").Value = "'1.080.43"
$ws.Range("E43").Value = "  -1.80%  "

$ws.Range("D44").Value = "'0.4272"
$ws.Range("E44").Value = "  -6.68%  "

$ws.Range("D45").Value = "'69.66"
$ws.Range("E45").Value = "  -4.05%  "

$ws.Range("D46").Value = "'0.9999"

$ws.Range("D47").Value = "'101.90"
$ws.Range("E47").Value = "  -2.68%  "

$ws.Range("E48").Value = "  -6.96%  "

$ws.Range("D49").Value = "'1.990.78"
$ws.Range("E49").Value = "  -5.50%  "

$ws.Range("D50").Value = "'7.019"
$ws.Range("E50").Value = "  -7.69%  "

$ws.Range("D51").Value = "'9.110"
$ws.Range("E51").Value = "  -5.71%  "
